$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.530.07"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.919.70"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  +0.80%  "
$ws.Range("D5").Value = "325.99"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("D7").Value = "0.4816"
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("D8").Value = "0.4059"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("D9").Value = "0.08224"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("D11").Value = "23.42"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "1.943.38"
$ws.Range("D13").Value = "6.059"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").Value = "7.246"
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("D15").Value = "91.54"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").Value = "0.06867"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").Value = "1.013"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "0.00001040"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "17.56"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").Value = "29.543.09"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").Value = "11.85"
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("D24").Value = "2.194"
$ws.Range("E24").Value = "  +1.35%  "
$ws.Range("D25").Value = "2.169.37"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("D26").Value = "6.558"
$ws.Range("E26").Value = "  +4.54%  "
$ws.Range("D27").Value = "156.06"
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("D28").Value = "20.04"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "2.100"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "120.63"
$ws.Range("E30").Value = "  +0.61%  "
$ws.Range("D31").Value = "1.019"
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("D32").Value = "0.09638"
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("D33").Value = "5.620"
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("D34").Value = "3.560"
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("D35").Value = "1.374"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("D36").Value = "0.06310"
$ws.Range("E36").Value = "  +3.26%  "
$ws.Range("D37").Value = "0.02287"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("D38").Value = "1.185"
$ws.Range("E38").Value = "  +1.20%  "
$ws.Range("D39").Value = "0.5943"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "10.73"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").Value = "7.907"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("D43").Value = "2.460"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "12.36"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "0.07469"
$ws.Range("E46").Value = "  -3.18%  "
$ws.Range("D47").Value = "0.5565"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").Value = "1.944"
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("D49").Value = "118.67"
$ws.Range("E49").Value = "  +3.36%  "
$ws.Range("D50").Value = "2.430"
$ws.Range("E50").Value = "  +3.44%  "
$ws.Range("D51").Value = "72.16"
$ws.Range("E51").Value = "  -0.67%  "

$priceRange.ClearFormats()
